$wb = $excel.ActiveWorkbook

# The last three source files (rows 5-7 in each sheet) have finished the handoff
# step and moved on to translation, so their "Status" changes from
# "Ready for handoff" to "In Translation" (matching the other rows already
# showing "In Translation").

# Overview sheet: status is split across the zh-cn (E) and de-de (F) columns.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E5:F7").Value = "In Translation"

# zh-cn sheet: Status is column C.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C5:C7").Value = "In Translation"

# de-de sheet: Status is column C.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C5:C7").Value = "In Translation"
